$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row of data for "Alone"
$ws.Range("A2").Value = "Alone"
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0

# Update selection as seen in diff
$ws.Range("F3").Select() | Out-Null
